# Adds three more attendance-date column pairs (05-04-2025, 09-04-2025,
# 10-04-2025) to the FY1 student attendance sheet, and records the
# 07-04-2025 check-in for the two students who actually attended that day.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells F1:K1 -------------------------------------------------
# Clone the formatting already used by the existing header cell (E1 —
# bold, bordered, centered) across the six new header cells, then fill
# in their text.
$ws.Range("E1").Copy($ws.Range("F1:K1"))
$ws.Range("F1").Value = "05-04-2025 Status"
$ws.Range("G1").Value = "05-04-2025 Time"
$ws.Range("H1").Value = "09-04-2025 Status"
$ws.Range("I1").Value = "09-04-2025 Time"
$ws.Range("J1").Value = "10-04-2025 Status"
$ws.Range("K1").Value = "10-04-2025 Time"

# --- New data cells F2:K24 ---------------------------------------------------
# Every student gets the same placeholder attendance data for the three
# new dates: absent on 05-04 (marked "Off" with a blank time), and
# absent ("A") with a zero time for 09-04 and 10-04.
for ($row = 2; $row -le 24; $row++) {
    $ws.Cells.Item($row, 6).Value = "Off"            # F - 05-04-2025 Status
    $ws.Cells.Item($row, 7).Value = "'"               # G - 05-04-2025 Time (blank text cell)
    $ws.Cells.Item($row, 7).ClearFormats()            # drop the quote-prefix formatting the trick above leaves behind
    $ws.Cells.Item($row, 8).Value = "A"               # H - 09-04-2025 Status
    $ws.Cells.Item($row, 9).Value = "00:00:00"        # I - 09-04-2025 Time
    $ws.Cells.Item($row, 10).Value = "A"              # J - 10-04-2025 Status
    $ws.Cells.Item($row, 11).Value = "00:00:00"       # K - 10-04-2025 Time
}

# --- Corrected 07-04-2025 attendance for rows 18 & 19 -----------------------
# These two students actually showed up and checked in that day.
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = "P"
$ws.Range("E18").Value = "08:46:24 AM"

$ws.Range("C19").Value = 1
$ws.Range("D19").Value = "P"
$ws.Range("E19").Value = "08:46:03 AM"
